$p = $ppt.ActivePresentation

# Add a new slide at the end, using the "Title Only" layout (slideLayout6.xml),
# matching the layout used by the other content slides in this deck.
$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 6)

# Title
$s.Shapes.Item(1).TextFrame.TextRange.Text = 'Endpoint → REFDB tables, columns, and sources'

# Table: 7 rows x 4 columns, sized/positioned to match the target layout.
# (EMU values from the target OOXML divided by 12700 EMU-per-point.)
$tbl = $s.Shapes.AddTable(7, 4, 21.6, 93.6, 928.8, 374.4)
$tblShape = $tbl.Table

$tblShape.Columns.Item(1).Width = 232.2
$tblShape.Columns.Item(2).Width = 232.2
$tblShape.Columns.Item(3).Width = 232.2
$tblShape.Columns.Item(4).Width = 232.2

$tblShape.Rows.Item(1).Height = 53.48566929133858
$tblShape.Rows.Item(2).Height = 53.48566929133858
$tblShape.Rows.Item(3).Height = 53.48566929133858
$tblShape.Rows.Item(4).Height = 53.48566929133858
$tblShape.Rows.Item(5).Height = 53.48566929133858
$tblShape.Rows.Item(6).Height = 53.48566929133858
$tblShape.Rows.Item(7).Height = 53.4859842519685

$rows = @(
    @('Endpoint', 'Final REFDB table(s)', 'Key columns populated', 'Primary sources'),
    @('OnLot', 'ON_LOT; ON_PROD', 'lot, mfgLot, product, fab, sourceLot, lotType, maskSet, process, technology, PTI, family', 'LotG (native+WS), LTM WS (lotType), Data Warehouse (PLM/MfgArea), MES (Torrent/Genesis)'),
    @('OnProd', 'ON_PROD', 'product, productVersion, family, process, technology, maskSet, fab', 'MES (Torrent/Genesis), Data Warehouse PLM/MfgArea, LotG'),
    @('OnScribe', 'ON_SCRIBE', 'lot, waferNum, waferId, scribeId, insertTime, status', 'VID↔SCRIBE services (fab-configured); calculated fallback via AttributeUtils; OnLot cache for sourceLot context'),
    @('OnSlice', 'ON_SLICE', 'slice, globalWaferId, puckId, runId, sliceSourceLot, startLot, fabWaferId, fabSourceLot, slicePartname, sliceLottype, sliceSupplierId, puckHeight, sliceOrder, sliceStartTime', 'Primary writes via admin DTO/API; upstream ingestion uses BIWMES+eCofA+TORRENT to populate/maintain rows'),
    @('OnWmap', 'ON_WMAP', 'idWaferMapConfiguration, product/device mapping, metadata per WMC', 'Matchup service (by lot/scribe) and WMC service (by config/product) via Caller'),
    @('PP_LOTPROD (context)', 'PP_LOTPROD', 'lot, product, fab (frontend provenance)', 'Internal PP_LOTPROD DB exposed via /api/pplotprod/bylotid; consumed by ingestion scripts')
)

for ($r = 1; $r -le 7; $r++) {
    $rowData = $rows[$r - 1]
    for ($c = 1; $c -le 4; $c++) {
        $cell = $tblShape.Cell($r, $c)
        $tr = $cell.Shape.TextFrame.TextRange
        $tr.Text = $rowData[$c - 1]
        if ($r -eq 1) {
            $tr.Font.Bold = -1
            $tr.Font.Size = 14
        } else {
            $tr.Font.Size = 12
        }
    }
}
